$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.238.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.373.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.49%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.504"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.36%  "

$ws.Range("E11").Value = "  +4.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0788"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.743.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.377.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.802"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.232.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0889"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.27%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.78%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -13.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.94%  "

$ws.Range("E33").Value = "  +3.50%  "

$ws.Range("E34").Value = "  +10.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0732"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "127.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.34%  "

$ws.Range("E40").Value = "  -2.61%  "

$ws.Range("E41").Value = "  -0.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.930.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.33%  "

$ws.Range("E44").Value = "  -0.82%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.41%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.76%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.51%  "

$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.599.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.12%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.22%  "

$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.22%  "
